$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 13, pushing existing rows 13-21 down to 14-22,
# and carrying the row 13 formatting (incl. the date style in column D) down.
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the weekly record that was added.
$ws.Cells.Item(13, 1).Value = 2
$ws.Cells.Item(13, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(13, 3).Value = "Coquimbo"
$ws.Cells.Item(13, 4).Value = 44413
$ws.Cells.Item(13, 5).Value = 4
$ws.Cells.Item(13, 6).Value = 100112026
$ws.Cells.Item(13, 7).Value = "Haba"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 1200
$ws.Cells.Item(13, 11).Value = 10000
$ws.Cells.Item(13, 12).Value = 11000
$ws.Cells.Item(13, 13).Value = 10500
$ws.Cells.Item(13, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(13, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(13, 16).Value = 420
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = "Hortaliza"
